# The bulleted list originally contained four items:
#   Zegar
#   Przycisk do wjazdu i schowania opcji
#   Opcje
#   Przycisk do potwierdzenia opcji
#
# The edit trims the list down to a single remaining item:
#   Przycisk do wjazdu i schowania opcji
#
# Remove the other three list paragraphs (including their paragraph
# marks) so the surrounding paragraphs merge/close up correctly, which
# mirrors selecting+deleting those bullet lines in Word.

$d = $word.ActiveDocument

$toRemove = @("Zegar", "Opcje", "Przycisk do potwierdzenia opcji")

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()
    if ($toRemove -contains $text) {
        $para.Range.Delete()
    }
}
